# 7.62x39 damage increased from 40 to 45
# Row 19 = ammo_7.62x39_fmj, Row 20 = ammo_7.62x39_ap
# Column H ("ammo_k_hit") drives the total damage (J) and in-game displayed
# power (K) via the existing formulas, so only H needs to be changed; the
# dependent formula cells recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("H19").Value = 1.17
$ws.Range("H20").Value = 1.17

# Update the saved view/selection state to match the authored workbook
# (scrolled back to the top, selection moved to K13).
$ws.Activate()
$ws.Range("K13").Select()
